# "Call BT func working"
# Re-measured the BT power numbers: the existing "3_3" sheet's Deep_Sleep /
# BT_Idle / BT_Iscan rows get the new averages + raw samples (and BT_Iscan's
# label is corrected to BT_I_Scan), and a second sheet "1_8" is added with
# the same table shape for the new 1.8 scan pass.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Sheet 1 "3_3": fix the row label, then refresh the numbers ----
$ws1.Range("A4").Value = "BT_I_Scan"

# Deep_Sleep row (row 2)
$ws1.Range("B2").Value = 0.2398
$ws1.Range("C2").Value = 0.2476
$ws1.Range("D2").Value = 0.2345
$ws1.Range("E2").Value = 0.0043
$ws1.Range("G2").Value = "0.238239,0.235705,0.24335,0.23686,0.245387,0.2374,0.247623,0.236842,0.234507,0.242493"

# BT_Idle row (row 3)
$ws1.Range("B3").Value = 8.3706
$ws1.Range("C3").Value = 8.387499999999999
$ws1.Range("D3").Value = 8.351800000000001
$ws1.Range("E3").Value = 0.0114
$ws1.Range("G3").Value = "8.361353,8.365427,8.383853,8.387523,8.366563,8.351752,8.370377,8.38453,8.375848,8.358701"

# BT_I_Scan row (row 4)
$ws1.Range("B4").Value = 0.2397
$ws1.Range("C4").Value = 0.2442
$ws1.Range("D4").Value = 0.231
$ws1.Range("E4").Value = 0.0039
$ws1.Range("G4").Value = "0.230954,0.23822,0.244188,0.243151,0.235724,0.240655,0.242294,0.24189,0.237978,0.242331"

# ---- Sheet 2 "1_8": new sheet, positioned right after "3_3" ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "1_8"

# Pull header row + row labels over with their existing formatting instead of
# re-declaring styles (keeps the same bold/centered/bordered look).
$ws1.Range("B1:G1").Copy($ws2.Range("B1:G1"))
$ws1.Range("A2:A4").Copy($ws2.Range("A2:A4"))

$ws2.Range("B2").Value = -0.1003
$ws2.Range("C2").Value = -0.08599999999999999
$ws2.Range("D2").Value = -0.1288
$ws2.Range("E2").Value = 0.0173
$ws2.Range("F2").Value = 10
$ws2.Range("G2").Value = "-0.091836,-0.088261,-0.086551,-0.127434,-0.097417,-0.087403,-0.122409,-0.086822,-0.086014,-0.128823"

$ws2.Range("B3").Value = -5.5091
$ws2.Range("C3").Value = -5.4969
$ws2.Range("D3").Value = -5.5191
$ws2.Range("E3").Value = 0.0073
$ws2.Range("F3").Value = 10
$ws2.Range("G3").Value = "-5.506857,-5.5191,-5.511889,-5.503017,-5.502147,-5.516316,-5.519082,-5.511321,-5.49688,-5.504369"

$ws2.Range("B4").Value = -0.5298
$ws2.Range("C4").Value = -0.0827
$ws2.Range("D4").Value = -4.4325
$ws2.Range("E4").Value = 1.301
$ws2.Range("F4").Value = 10
$ws2.Range("G4").Value = "-4.43255,-0.082661,-0.094892,-0.124082,-0.088526,-0.084347,-0.125175,-0.086841,-0.087958,-0.091298"

# Keep "3_3" the selected/active tab, matching the source workbook.
$ws1.Activate()
